$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update Correspond Handoff Datetime (D2) and
# Correspond Handback DateTime (G2) for the 1c9b1662... row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-11 03:44:53"
$wsZhCn.Range("G2").Value = "2016-01-11 03:46:17"

# "de-de" sheet: update Correspond Handoff Datetime (D2) and
# Correspond Handback DateTime (G2) for the 1c9b1662... row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-11 03:45:17"
$wsDeDe.Range("G2").Value = "2016-01-11 03:46:54"
